# edit.ps1 -- applies the "architecture_platform" edits described by the
# target diff:
#
#  1. In "We choose Node as our architecture/platform for developing
#     solution for use cases. Since we choose Javascript..." the word
#     "Since " is deleted and the following "we" is capitalized to "We"
#     (so the sentence reads "...for use cases. We choose Javascript...").
#     The target XML shows this run split into three runs (the prefix
#     sentence, a lone "W", and "e choose ") rather than one fused run,
#     so we perform the edit as two ops: delete "Since ", then rewrite
#     just the leading "w" -> "W" with a formatting no-op (Bold on/off)
#     which forces the engine to keep it as a distinct run instead of
#     silently re-merging it with its identically-formatted neighbours.
#
#  2-4. Three spots where several adjacent, identically formatted runs
#     get fused into a single run with the same overall text (no visible
#     text change, just run consolidation). We reproduce this by doing a
#     Find/Replace of the full span with itself, which causes the engine
#     to rewrite that span as one run.

$d = $word.ActiveDocument

# --- Change 1a: remove "Since " -------------------------------------
$ok1 = $d.Content.Find.Execute("Since ", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

# --- Change 1b: capitalize "we choose" -> "We choose", as its own run -
$rngWe = $d.Content
$ok2 = $rngWe.Find.Execute("we choose", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$firstChar = $d.Range($rngWe.Start, $rngWe.Start + 1)
$firstChar.Text = "W"
# toggle formatting to force this single character to stay its own run
# rather than being silently re-merged into its neighbours
$firstChar.Bold = 1
$firstChar.Bold = 0

# --- Change 2: consolidate " and back end -- look at slide show from
#     10/31/16 for reason why MEAN stack is good" into a single run ----
$t2 = " and back end -- look at slide show from 10/31/16 for reason why MEAN stack is good"
$ok3 = $d.Content.Find.Execute($t2, $true, $false, $false, $false, $false, $true, 1, $false, $t2, 2)

# --- Change 3: consolidate the PHP paragraph's opening sentence run --
$rsq = [char]0x2019
$t3 = "We are not choosing PHP. First, we can" + $rsq + "t work on front end and back end independently as we do in Node. Second, PHP is open source, which is insecure since it pays little attention to security. This makes it inappropriate for our project especially we have to deal with both YouTube and "
$ok4 = $d.Content.Find.Execute($t3, $true, $false, $false, $false, $false, $true, 1, $false, $t3, 2)

# --- Change 4: consolidate the trailing "accounts and their passwords;
#     we don't want..." run -------------------------------------------
$t4 = " accounts and their passwords; we don" + $rsq + "t want any of this sensitive account information to land in the wrong hands."
$ok5 = $d.Content.Find.Execute($t4, $true, $false, $false, $false, $false, $true, 1, $false, $t4, 2)

Write-Output "Since-removed: $ok1; we-found: $ok2; merge2: $ok3; merge3: $ok4; merge4: $ok5"
